$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 123, shifting existing rows 123:234 down to 124:235
$ws.Rows.Item(123).Insert()

# Populate the newly inserted row 123 with the new price record
$ws.Range("A123").Value = 7
$ws.Range("B123").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C123").Value = "Ñuble"
$ws.Range("D123").Value = 45049
$ws.Range("E123").Value = 16
$ws.Range("F123").Value = "Fruta"
$ws.Range("G123").Value = 100102
$ws.Range("H123").Value = "Cítricos"
$ws.Range("I123").Value = 100102004
$ws.Range("J123").Value = "Mandarina"
$ws.Range("K123").Value = "Murcott"
$ws.Range("L123").Value = "Primera"
$ws.Range("M123").Value = 120
$ws.Range("N123").Value = 14000
$ws.Range("O123").Value = 15000
$ws.Range("P123").Value = 14500
$ws.Range("Q123").Value = "$/caja 12 kilos"
$ws.Range("R123").Value = "Región de O'Higgins"
$ws.Range("S123").Value = 1208
$ws.Range("T123").Value = 12
